$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the displayed period values in E16:E22 (2308..2302 -> 2302..2308)
$ws.Range("E16").Value = "2302"
$ws.Range("E17").Value = "2303"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2305"
$ws.Range("E20").Value = "2306"
$ws.Range("E21").Value = "2307"
$ws.Range("E22").Value = "2308"

# Swap the date values between F16 and F22
$ws.Range("F16").Value = 44854
$ws.Range("F22").Value = 46400
